# ADD results from server
# Update the computed investment-capacity results for years 2025, 2030 and 2035
# with the latest values returned from the server.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.003837439598571024
$ws.Range("E2").Value = 0.3718167378372526
$ws.Range("G2").Value = 0.2494892361375017
$ws.Range("I2").Value = 0.3687520809646043
$ws.Range("L2").Value = 0.597153
$ws.Range("M2").Value = 0.0822565
$ws.Range("N2").Value = 12.82009457445582
$ws.Range("O2").Value = 3.538061458581253

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0.009260380363639123
$ws.Range("B2").Value = 0.04893136040142915
$ws.Range("E2").Value = 0.2216906708742613
$ws.Range("I2").Value = 0.4247748743219193
$ws.Range("L2").Value = 0.1116199591040388
$ws.Range("M2").Value = 0.04737166666666669
$ws.Range("N2").Value = 5.019473591562843
$ws.Range("O2").Value = 2.349307442578966

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.08320519174800101
$ws.Range("B2").Value = 0.02872957743582358
$ws.Range("E2").Value = 0.1707961561752257
$ws.Range("I2").Value = 0.4661253951561688
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04858258333333332
$ws.Range("N2").Value = 8.420331737781321
$ws.Range("O2").Value = 4.913815709442911
